$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Row 73: new entry dated 13.1.2022 (serial 44574), 1h worked, comment, project "client"
$ws.Cells.Item(73, 1).Value = 44574
$ws.Cells.Item(73, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(73, 2).Value = 1
$ws.Cells.Item(73, 3).Value = "kommenttien lisäystä melkein kaikkiin komponentteihin, myös vääriä kommentteja korjattu, pientä refaktorointia"
$ws.Cells.Item(73, 4).Value = "client"

# Row 74: second new entry, 1h worked, comment, project "client"
$ws.Cells.Item(74, 2).Value = 1
$ws.Cells.Item(74, 3).Value = "/home pohja, Info ja TopCoins aloitettu"
$ws.Cells.Item(74, 4).Value = "client"

# Update the total row formula to include the new rows
$ws.Cells.Item(75, 2).Formula = "=SUM(B2:B74)"

# Move the active selection to C73, matching the new last data entry
$ws.Range("C73").Select()
